$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "03/23/2020 22:38:38"
